$d = $word.ActiveDocument
$nl = [char]10

$lines = @(
    "*",
    "UNIVERSTTY OF ABUJA,",
    "DEPARTMENT: EDUCATIONAL FOUNDATIONS,",
    "DESCRIPTION: CURRICULUM AND INSTRUCTION,",
    "COURSE CODE: EDU 205,",
    "CREDIT UNIT: 2",
    "NSTRUCTION: ANSWER ANY THREE QUESTIONS,",
    "LEVEL: 200,",
    "FERST SEMESTER EXAMINATION 2018/2019",
    "TIME ALLOWED: 2 BHOURS.",
    "a. Define curiculum and show how it is related to instruction",
    "1b. State five characteristics of lesson pla",
    "2a. Explain the meaning of curriculum design",
    "2b. Write short notes on the following types of curriculum designs",
    "a) Subject centred",
    "b) Activity centred",
    "CCore- centred",
    "Ba. Discuss four factors affecting effective teachimg",
    "3b. Identify and discuss four ypes of motivation",
    "4a. Explain the term classroom management",
    "45. Briefly discuSS four roies of teachers in Ciassroom management",
    "5. State and explain Benjamin Bloom's taxonomy of educational objectives"
)

$newText = [string]::Join($nl, $lines)

$d.Paragraphs(1).Range.Text = $newText
